$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "39.932.09"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -0.05%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.233.32"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -4.59%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.16%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "291.83"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -5.71%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "85.56"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +1.58%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.515"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -1.82%  "

# Row 8
$ws.Range("E8").Value = "  -0.18%  "

# Row 9
$ws.Range("E9").Value = "  -1.94%  "

# Row 10
$ws.Range("E10").Value = "  -0.53%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "30.43"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +1.81%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "47.82"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -8.89%  "

# Row 13
$ws.Range("E13").Value = "  -2.31%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.37"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -0.37%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.574.05"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -4.90%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.20"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -3.75%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.220.80"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -5.81%  "

# Row 18
$ws.Range("E18").Value = "  -4.12%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "39.820.07"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.40%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0892"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -0.71%  "

# Row 21
$ws.Range("E21").Value = "  -4.28%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.60"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.65%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "65.47"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -3.57%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "232.74"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -0.87%  "

# Row 25
$ws.Range("E25").Value = "  +0.08%  "

# Row 26
$ws.Range("E26").Value = "  -4.82%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.85"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +1.54%  "

# Row 28
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "22.98"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -2.13%  "

# Row 29
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.13"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -2.98%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.26"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +0.56%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "155.07"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +1.22%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "33.21"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -4.02%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.999"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -0.36%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.85"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -4.29%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0709"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -0.71%  "

# Row 36
$ws.Range("E36").Value = "  -4.91%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "16.56"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +6.96%  "

# Row 38
$ws.Range("E38").Value = "  -1.10%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0980"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -0.53%  "

# Row 40
$ws.Range("E40").Value = "  -3.35%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.67"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -2.55%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.75"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -2.20%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.955.37"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -0.83%  "

# Row 44
$ws.Range("E44").Value = "  -2.98%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0269"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +2.08%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "9.52"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +1.32%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "16.24"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -7.18%  "

# Row 48
$ws.Range("E48").Value = "  -2.31%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.446.02"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -4.66%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "71.02"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +1.56%  "

# Row 51
$ws.Range("E51").Value = "  +7.17%  "
